$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 335.57144
$ws.Range("I5").Value = 149.8
$ws.Range("K5").Value = 149.8
$ws.Range("M5").Value = -34.80000000000001
$ws.Range("H12").Value = 551.4286
$ws.Range("I12").Value = 486.66666
$ws.Range("J12").Value = 600
$ws.Range("K12").Value = 486.66666
$ws.Range("L12").Value = 600
$ws.Range("M12").Value = -316.66666
$ws.Range("N12").Value = -940
$ws.Range("H32").Value = 2454.6
$ws.Range("I32").Value = 4148.6665
$ws.Range("J32").Value = 1728.5714
$ws.Range("K32").Value = 4148.6665
$ws.Range("L32").Value = 1728.5714
$ws.Range("M32").Value = -3822.6665
$ws.Range("N32").Value = -2380.5714
$ws.Range("H74").Value = 4360
$ws.Range("I74").Value = 800
$ws.Range("J74").Value = 5250
$ws.Range("K74").Value = 800
$ws.Range("L74").Value = 5250
$ws.Range("M74").Value = 136
$ws.Range("N74").Value = -7122
$ws.Range("H77").Value = 4360
$ws.Range("I77").Value = 800
$ws.Range("J77").Value = 5250
$ws.Range("K77").Value = 4000
$ws.Range("L77").Value = 26250
$ws.Range("M77").Value = 680
$ws.Range("N77").Value = -35610
$ws.Range("H132").Value = 1483506
$ws.Range("I132").Value = 2021406.4
$ws.Range("J132").Value = 4280.125
$ws.Range("K132").Value = 6064219.199999999
$ws.Range("L132").Value = 12840.375
$ws.Range("M132").Value = -6061689.199999999
$ws.Range("N132").Value = -17900.375
$ws.Range("H135").Value = 1222.8
$ws.Range("I135").Value = 1025.4286
$ws.Range("K135").Value = 9228.857399999999
$ws.Range("M135").Value = -6693.857399999999
$ws.Range("H137").Value = 57360.89
$ws.Range("I137").Value = 1949.7693
$ws.Range("K137").Value = 5849.3079
$ws.Range("M137").Value = -3299.3079
$ws.Range("H141").Value = 52546.332
$ws.Range("I141").Value = 62315.332
$ws.Range("K141").Value = 186945.996
$ws.Range("M141").Value = -181765.996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1237.1177
$ws.Range("I97").Value = 1005.2857
$ws.Range("J97").Value = 2319
$ws.Range("K97").Value = 1005.2857
$ws.Range("L97").Value = 2319
$ws.Range("M97").Value = -509.2857
$ws.Range("N97").Value = -3311

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7362.04
$ws.Range("I20").Value = 5150.3687
$ws.Range("K20").Value = 5150.3687
$ws.Range("M20").Value = -4903.3687
$ws.Range("H86").Value = 2012.909
$ws.Range("I86").Value = 1940.8334
$ws.Range("J86").Value = 2337.25
$ws.Range("K86").Value = 1940.8334
$ws.Range("L86").Value = 2337.25
$ws.Range("M86").Value = -817.8334
$ws.Range("N86").Value = -4583.25
$ws.Range("H89").Value = 2012.909
$ws.Range("I89").Value = 1940.8334
$ws.Range("J89").Value = 2337.25
$ws.Range("K89").Value = 9704.166999999999
$ws.Range("L89").Value = 11686.25
$ws.Range("M89").Value = -4088.166999999999
$ws.Range("N89").Value = -22918.25
$ws.Range("H94").Value = 1150.6957
$ws.Range("I94").Value = 1165.8948
$ws.Range("K94").Value = 1165.8948
$ws.Range("M94").Value = -714.8948

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 460.86365
$ws.Range("I7").Value = 341.42856
$ws.Range("K7").Value = 341.42856
$ws.Range("M7").Value = -228.42856
$ws.Range("H31").Value = 2125
$ws.Range("I31").Value = 1824.4286
$ws.Range("K31").Value = 1824.4286
$ws.Range("M31").Value = -1529.4286
$ws.Range("H34").Value = 2125
$ws.Range("I34").Value = 1824.4286
$ws.Range("K34").Value = 1824.4286
$ws.Range("M34").Value = -1622.4286
$ws.Range("H134").Value = 2334.5715
$ws.Range("I134").Value = 1865.6
$ws.Range("K134").Value = 5596.799999999999
$ws.Range("M134").Value = -3061.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 97.36364
$ws.Range("I2").Value = 82.09999999999999
$ws.Range("K2").Value = 492.6
$ws.Range("M2").Value = -379.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 324.375
$ws.Range("I2").Value = 282.33334
$ws.Range("J2").Value = 450.5
$ws.Range("K2").Value = 282.33334
$ws.Range("L2").Value = 450.5
$ws.Range("M2").Value = -169.33334
$ws.Range("N2").Value = -676.5
$ws.Range("H32").Value = 70644.5
$ws.Range("J32").Value = 70644.5
$ws.Range("L32").Value = 70644.5
$ws.Range("N32").Value = -71236.5
$ws.Range("H42").Value = 40000
$ws.Range("J42").Value = 40000
$ws.Range("L42").Value = 40000
$ws.Range("N42").Value = -40970
$ws.Range("H80").Value = 8216.174000000001
$ws.Range("I80").Value = 9612.666999999999
$ws.Range("K80").Value = 9612.666999999999
$ws.Range("M80").Value = -8614.666999999999
$ws.Range("H83").Value = 8216.174000000001
$ws.Range("I83").Value = 9612.666999999999
$ws.Range("K83").Value = 48063.335
$ws.Range("M83").Value = -43071.335
$ws.Range("H115").Value = 40000
$ws.Range("J115").Value = 40000
$ws.Range("L115").Value = 40000
$ws.Range("N115").Value = -42350
$ws.Range("H126").Value = 3602.8
$ws.Range("I126").Value = 3333.3333
$ws.Range("K126").Value = 9999.999899999999
$ws.Range("M126").Value = -7529.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3315.75
$ws.Range("I16").Value = 932
$ws.Range("K16").Value = 932
$ws.Range("M16").Value = -762
$ws.Range("H22").Value = 1535.25
$ws.Range("I22").Value = 1535.25
$ws.Range("K22").Value = 1535.25
$ws.Range("M22").Value = -1240.25
$ws.Range("H27").Value = 1535.25
$ws.Range("I27").Value = 1535.25
$ws.Range("K27").Value = 1535.25
$ws.Range("M27").Value = -1428.25
$ws.Range("H100").Value = 96435.42999999999
$ws.Range("I100").Value = 1676.6
$ws.Range("K100").Value = 1676.6
$ws.Range("M100").Value = -1135.6
$ws.Range("H122").Value = 3797.158
$ws.Range("I122").Value = 3602.875
$ws.Range("K122").Value = 10808.625
$ws.Range("M122").Value = -8358.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1650.6666
$ws.Range("I96").Value = 1474
$ws.Range("J96").Value = 2004
$ws.Range("K96").Value = 1474
$ws.Range("L96").Value = 2004
$ws.Range("M96").Value = -101
$ws.Range("N96").Value = -4750
$ws.Range("H132").Value = 7660.2256
$ws.Range("I132").Value = 9162.137000000001
$ws.Range("K132").Value = 27486.411
$ws.Range("M132").Value = -24956.411
